$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 801
$ws.Range("D2").Value = 1063
$ws.Range("E2").Value = 1864
$ws.Range("F2").Value = 1144.38307919425
$ws.Range("G2").Value = -678.0789917484562
$ws.Range("H2").Value = 0.1162664871354343
$ws.Range("I2").Value = -0.0682687876290619
$ws.Range("J2").Value = 0.01102990072199039
$ws.Range("K2").Value = 916650.8464345935
$ws.Range("L2").Value = -720797.9682286098
$ws.Range("M2").Value = 0.4297210300429185
$ws.Range("N2").Value = 1.687684020770808
$ws.Range("O2").Value = 1.271716745660785
$ws.Range("P2").Value = 195852.8782059838

$ws.Range("C3").Value = 880
$ws.Range("D3").Value = 1367
$ws.Range("E3").Value = 2247
$ws.Range("F3").Value = 1033.914640218658
$ws.Range("G3").Value = -676.5741857166947
$ws.Range("H3").Value = 0.1097089115712895
$ws.Range("I3").Value = -0.06786639221842332
$ws.Range("J3").Value = 0.001678008019648457
$ws.Range("K3").Value = 909844.8833924202
$ws.Range("L3").Value = -924876.9118747229
$ws.Range("M3").Value = 0.3916332888295505
$ws.Range("N3").Value = 1.528161526179768
$ws.Range("O3").Value = 0.983746995638768
$ws.Range("P3").Value = -15032.02848230267

$ws.Range("C4").Value = 1534
$ws.Range("D4").Value = 1827
$ws.Range("E4").Value = 3361
$ws.Range("F4").Value = 1598.114591890394
$ws.Range("G4").Value = -644.8157321045246
$ws.Range("H4").Value = 0.1708592716548004
$ws.Range("I4").Value = -0.06377776658843408
$ws.Range("J4").Value = 0.04331334220809127
$ws.Range("K4").Value = 2451507.783959871
$ws.Range("L4").Value = -1178078.342554965
$ws.Range("M4").Value = 0.4564117822076763
$ws.Range("N4").Value = 2.478405088961663
$ws.Range("O4").Value = 2.080937825105202
$ws.Range("P4").Value = 1273429.441404906

$ws.Range("C5").Value = 1033
$ws.Range("D5").Value = 1612
$ws.Range("E5").Value = 2645
$ws.Range("F5").Value = 1267.20155353264
$ws.Range("G5").Value = -638.0968716846857
$ws.Range("H5").Value = 0.136836432055763
$ws.Range("I5").Value = -0.0621039458797431
$ws.Range("J5").Value = 0.0155918614576399
$ws.Range("K5").Value = 1309019.204799217
$ws.Range("L5").Value = -1028612.157155714
$ws.Range("M5").Value = 0.3905482041587902
$ws.Range("N5").Value = 1.985907798273653
$ws.Range("O5").Value = 1.272607168496701
$ws.Range("P5").Value = 280407.0476435027

$ws.Range("C6").Value = 940
$ws.Range("D6").Value = 1458
$ws.Range("E6").Value = 2398
$ws.Range("F6").Value = 1155.481502181168
$ws.Range("G6").Value = -574.5511275784007
$ws.Range("H6").Value = 0.1173086308843675
$ws.Range("I6").Value = -0.05740118863322429
$ws.Range("J6").Value = 0.0110838949141219
$ws.Range("K6").Value = 1086152.612050297
$ws.Range("L6").Value = -837695.5440093079
$ws.Range("M6").Value = 0.3919933277731443
$ws.Range("N6").Value = 2.011103010190327
$ws.Range("O6").Value = 1.296595905060979
$ws.Range("P6").Value = 248457.0680409896

$ws.Range("C7").Value = 1044
$ws.Range("D7").Value = 1742
$ws.Range("E7").Value = 2786
$ws.Range("F7").Value = 1046.581711919139
$ws.Range("G7").Value = -528.6141877759942
$ws.Range("H7").Value = 0.1064986693698909
$ws.Range("I7").Value = -0.05288942590220822
$ws.Range("J7").Value = 0.006838202046130525
$ws.Range("K7").Value = 1092631.30724358
$ws.Range("L7").Value = -920845.9151057814
$ws.Range("M7").Value = 0.3747307968413496
$ws.Range("N7").Value = 1.97985929269579
$ws.Range("O7").Value = 1.18655172306223
$ws.Range("P7").Value = 171785.3921377991

$ws.Range("C8").Value = 786
$ws.Range("D8").Value = 1529
$ws.Range("E8").Value = 2315
$ws.Range("F8").Value = 1035.132609135504
$ws.Range("G8").Value = -569.9849725210119
$ws.Range("H8").Value = 0.1045555122542369
$ws.Range("I8").Value = -0.05722148318405369
$ws.Range("J8").Value = -0.002294175013644831
$ws.Range("K8").Value = 813614.2307805065
$ws.Range("L8").Value = -871507.0229846266
$ws.Range("M8").Value = 0.339524838012959
$ws.Range("N8").Value = 1.816070000156619
$ws.Range("O8").Value = 0.9335716285958818
$ws.Range("P8").Value = -57892.79220412008

$ws.Range("C9").Value = 698
$ws.Range("D9").Value = 1712
$ws.Range("E9").Value = 2410
$ws.Range("F9").Value = 1008.701179353257
$ws.Range("G9").Value = -779.1881600577467
$ws.Range("H9").Value = 0.103873949131873
$ws.Range("I9").Value = -0.07787987431049065
$ws.Range("J9").Value = -0.02523914038403022
$ws.Range("K9").Value = 704073.4231885733
$ws.Range("L9").Value = -1333970.130018862
$ws.Range("M9").Value = 0.2896265560165975
$ws.Range("N9").Value = 1.294554038498867
$ws.Range("O9").Value = 0.5278029899954492
$ws.Range("P9").Value = -629896.7068302884

$ws.Range("C10").Value = 1732
$ws.Range("D10").Value = 2599
$ws.Range("E10").Value = 4331
$ws.Range("F10").Value = 2072.291467075951
$ws.Range("G10").Value = -786.3332677659623
$ws.Range("H10").Value = 0.2108185888499272
$ws.Range("I10").Value = -0.07774070692665151
$ws.Range("J10").Value = 0.03765636079097364
$ws.Range("K10").Value = 3589208.82097555
$ws.Range("L10").Value = -2043680.162923735
$ws.Range("M10").Value = 0.3999076425767721
$ws.Range("N10").Value = 2.635385722600166
$ws.Range("O10").Value = 1.756247815137935
$ws.Range("P10").Value = 1545528.658051815

$ws.Range("C11").Value = 1654
$ws.Range("D11").Value = 2275
$ws.Range("E11").Value = 3929
$ws.Range("F11").Value = 1227.153118326884
$ws.Range("G11").Value = -580.3485045443944
$ws.Range("H11").Value = 0.1240811905108539
$ws.Range("I11").Value = -0.05782848206939645
$ws.Range("J11").Value = 0.01875044347087688
$ws.Range("K11").Value = 2029711.257712666
$ws.Range("L11").Value = -1320292.847838499
$ws.Range("M11").Value = 0.4209722575719013
$ws.Range("N11").Value = 2.114510692657453
$ws.Range("O11").Value = 1.537318982705679
$ws.Range("P11").Value = 709418.4098741664

$ws.Range("C12").Value = 1085
$ws.Range("D12").Value = 2030
$ws.Range("E12").Value = 3115
$ws.Range("F12").Value = 1084.791035308492
$ws.Range("G12").Value = -613.7745384458839
$ws.Range("H12").Value = 0.1160233992882252
$ws.Range("I12").Value = -0.06125215637653859
$ws.Range("J12").Value = 0.0004955090797273362
$ws.Range("K12").Value = 1176998.273309714
$ws.Range("L12").Value = -1245962.313045143
$ws.Range("M12").Value = 0.3483146067415731
$ws.Range("N12").Value = 1.767409638815015
$ws.Range("O12").Value = 0.9446499793666465
$ws.Range("P12").Value = -68964.03973542945

$ws.Range("C13").Value = 1462
$ws.Range("D13").Value = 2264
$ws.Range("E13").Value = 3726
$ws.Range("F13").Value = 1091.591217233229
$ws.Range("G13").Value = -555.0845890207484
$ws.Range("H13").Value = 0.1146091934118741
$ws.Range("I13").Value = -0.05543524584678286
$ws.Range("J13").Value = 0.01128643160790209
$ws.Range("K13").Value = 1595906.359594978
$ws.Range("L13").Value = -1256711.509542974
$ws.Range("M13").Value = 0.3923778851315083
$ws.Range("N13").Value = 1.966531297795454
$ws.Range("O13").Value = 1.269906694954484
$ws.Range("P13").Value = 339194.8500520042

$ws.Range("C14").Value = 1843
$ws.Range("D14").Value = 2452
$ws.Range("E14").Value = 4295
$ws.Range("F14").Value = 1065.884872302433
$ws.Range("G14").Value = -484.674831650149
$ws.Range("H14").Value = 0.1071775333938199
$ws.Range("I14").Value = -0.04815722206655487
$ws.Range("J14").Value = 0.01849748208093547
$ws.Range("K14").Value = 1964425.819653384
$ws.Range("L14").Value = -1188422.687206169
$ws.Range("M14").Value = 0.4291036088474971
$ws.Range("N14").Value = 2.199175205102906
$ws.Range("O14").Value = 1.652968965336316
$ws.Range("P14").Value = 776003.1324472157

$ws.Range("C15").Value = 1172
$ws.Range("D15").Value = 2070
$ws.Range("E15").Value = 3242
$ws.Range("F15").Value = 1040.768102913281
$ws.Range("G15").Value = -549.217332533061
$ws.Range("H15").Value = 0.1046125161859281
$ws.Range("I15").Value = -0.05457258570722654
$ws.Range("J15").Value = 0.002973663342365449
$ws.Range("K15").Value = 1219780.216614366
$ws.Range("L15").Value = -1136879.878343437
$ws.Range("M15").Value = 0.3615052436767427
$ws.Range("N15").Value = 1.895002290100942
$ws.Range("O15").Value = 1.072919171013673
$ws.Range("P15").Value = 82900.33827092894

$ws.Range("C16").Value = 1106
$ws.Range("D16").Value = 2311
$ws.Range("E16").Value = 3417
$ws.Range("F16").Value = 908.1881843908076
$ws.Range("G16").Value = -591.4227875012275
$ws.Range("H16").Value = 0.09146085340983208
$ws.Range("I16").Value = -0.058853642517681
$ws.Range("J16").Value = -0.01020048697310124
$ws.Range("K16").Value = 1004456.131936234
$ws.Range("L16").Value = -1366778.061915339
$ws.Range("M16").Value = 0.3236757389522973
$ws.Range("N16").Value = 1.53559890417466
$ws.Range("O16").Value = 0.7349080000074307
$ws.Range("P16").Value = -362321.9299791048

$ws.Range("C17").Value = 1992
$ws.Range("D17").Value = 2827
$ws.Range("E17").Value = 4819
$ws.Range("F17").Value = 1070.111489358152
$ws.Range("G17").Value = -547.4866193677728
$ws.Range("H17").Value = 0.1076160989662844
$ws.Range("I17").Value = -0.05439794744159979
$ws.Range("J17").Value = 0.01257278931799868
$ws.Range("K17").Value = 2131662.086801437
$ws.Range("L17").Value = -1547744.672952696
$ws.Range("M17").Value = 0.413363768416684
$ws.Range("N17").Value = 1.954589302280842
$ws.Range("O17").Value = 1.377269858557987
$ws.Range("P17").Value = 583917.4138487414

$ws.Range("C18").Value = 1675
$ws.Range("D18").Value = 2736
$ws.Range("E18").Value = 4411
$ws.Range("F18").Value = 1144.315954102867
$ws.Range("G18").Value = -535.7679475061136
$ws.Range("H18").Value = 0.1145522627810754
$ws.Range("I18").Value = -0.05318711397743991
$ws.Range("J18").Value = 0.01050897672093083
$ws.Range("K18").Value = 1916729.223122303
$ws.Range("L18").Value = -1465861.104376727
$ws.Range("M18").Value = 0.3797324869644071
$ws.Range("N18").Value = 2.135842503138562
$ws.Range("O18").Value = 1.307579017820574
$ws.Range("P18").Value = 450868.1187455754

$ws.Range("C19").Value = 1425
$ws.Range("D19").Value = 2640
$ws.Range("E19").Value = 4065
$ws.Range("F19").Value = 1200.111379838864
$ws.Range("G19").Value = -572.865577361677
$ws.Range("H19").Value = 0.1186258262530552
$ws.Range("I19").Value = -0.05703796611577005
$ws.Range("J19").Value = 0.004541592094703748
$ws.Range("K19").Value = 1710158.716270382
$ws.Range("L19").Value = -1512365.124234831
$ws.Range("M19").Value = 0.3505535055350554
$ws.Range("N19").Value = 2.094926676107784
$ws.Range("O19").Value = 1.130784285399085
$ws.Range("P19").Value = 197793.5920355509

$ws.Range("C20").Value = 1919
$ws.Range("D20").Value = 3314
$ws.Range("E20").Value = 5233
$ws.Range("F20").Value = 1145.732746977254
$ws.Range("G20").Value = -609.6339503149388
$ws.Range("H20").Value = 0.1141038991626735
$ws.Range("I20").Value = -0.06071198893315322
$ws.Range("J20").Value = 0.003394964870762576
$ws.Range("K20").Value = 2198661.141449352
$ws.Range("L20").Value = -2020326.911343707
$ws.Range("M20").Value = 0.3667112554939805
$ws.Range("N20").Value = 1.879378184868747
$ws.Range("O20").Value = 1.088269986953268
$ws.Range("P20").Value = 178334.2301056448

$ws.Range("C21").Value = 2508
$ws.Range("D21").Value = 3965
$ws.Range("E21").Value = 6473
$ws.Range("F21").Value = 1786.309786489659
$ws.Range("G21").Value = -784.473412366504
$ws.Range("H21").Value = 0.1772687735964914
$ws.Range("I21").Value = -0.0765433480535273
$ws.Range("J21").Value = 0.02179757595361754
$ws.Range("K21").Value = 4480064.944516069
$ws.Range("L21").Value = -3110437.080033188
$ws.Range("M21").Value = 0.3874555847365982
$ws.Range("N21").Value = 2.277081362261771
$ws.Range("O21").Value = 1.440332927251583
$ws.Range("P21").Value = 1369627.864482882

